# Dev IV Solo Project Rubric - mark directional & point light rows complete
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows that need to be marked Milestone "I" complete with "X"
# (cells in columns E & F already carry locked="0" protection style,
# so they can be edited while the sheet stays protected)
$rows = @(5, 7, 8, 9, 57, 66)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "I"
    $ws.Range("F$r").Value = "X"
}

# Restore the view: scroll back to top-left and select F9
$ws.Range("F9").Select()

$wb.Save()
